$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range('D2') '29.249.57'
Set-TextValue $ws.Range('E2') '  -0.29%  '
Set-TextValue $ws.Range('D3') '1.890.14'
Set-TextValue $ws.Range('E3') '  -0.89%  '
Set-TextValue $ws.Range('E4') '  +0.28%  '
Set-TextValue $ws.Range('D5') '322.71'
Set-TextValue $ws.Range('E5') '  -3.34%  '
Set-TextValue $ws.Range('D6') '1.001'
Set-TextValue $ws.Range('E6') '  +0.22%  '
Set-TextValue $ws.Range('D7') '0.4736'
Set-TextValue $ws.Range('E7') '  +2.11%  '
Set-TextValue $ws.Range('D8') '0.4039'
Set-TextValue $ws.Range('E8') '  -2.52%  '
Set-TextValue $ws.Range('D9') '47.27'
Set-TextValue $ws.Range('E9') '  -1.38%  '
Set-TextValue $ws.Range('D10') '0.08004'
Set-TextValue $ws.Range('E10') '  -0.53%  '
Set-TextValue $ws.Range('D11') '0.9965'
Set-TextValue $ws.Range('E11') '  -2.51%  '
Set-TextValue $ws.Range('D12') '23.28'
Set-TextValue $ws.Range('E12') '  +4.46%  '
Set-TextValue $ws.Range('D13') '1.893.48'
Set-TextValue $ws.Range('E13') '  -0.92%  '
Set-TextValue $ws.Range('D14') '5.931'
Set-TextValue $ws.Range('E14') '  -0.63%  '
Set-TextValue $ws.Range('D15') '7.027'
Set-TextValue $ws.Range('E15') '  -1.67%  '
Set-TextValue $ws.Range('D16') '89.18'
Set-TextValue $ws.Range('E16') '  -0.10%  '
Set-TextValue $ws.Range('D17') '1.001'
Set-TextValue $ws.Range('E17') '  +0.20%  '
Set-TextValue $ws.Range('D18') '0.06628'
Set-TextValue $ws.Range('E18') '  +0.60%  '
Set-TextValue $ws.Range('E19') '  -0.47%  '
Set-TextValue $ws.Range('E20') '  -0.94%  '
Set-TextValue $ws.Range('E21') '  +0.06%  '
Set-TextValue $ws.Range('D22') '29.277.15'
Set-TextValue $ws.Range('E22') '  -0.21%  '
Set-TextValue $ws.Range('E23') '  -0.57%  '
Set-TextValue $ws.Range('D24') '11.72'
Set-TextValue $ws.Range('E24') '  +2.53%  '
Set-TextValue $ws.Range('D25') '2.173'
Set-TextValue $ws.Range('E25') '  -1.12%  '
Set-TextValue $ws.Range('D26') '2.115.07'
Set-TextValue $ws.Range('E26') '  -3.56%  '
Set-TextValue $ws.Range('D27') '154.75'
Set-TextValue $ws.Range('E27') '  -1.24%  '
Set-TextValue $ws.Range('E28') '  -1.14%  '
Set-TextValue $ws.Range('D29') '5.933'
Set-TextValue $ws.Range('E29') '  +4.62%  '
Set-TextValue $ws.Range('E30') '  -3.11%  '
Set-TextValue $ws.Range('D31') '117.45'
Set-TextValue $ws.Range('E31') '  -0.04%  '
Set-TextValue $ws.Range('E32') '  -1.88%  '
Set-TextValue $ws.Range('D33') '0.09426'
Set-TextValue $ws.Range('E33') '  -0.39%  '
Set-TextValue $ws.Range('D34') '3.527'
Set-TextValue $ws.Range('E34') '  -0.52%  '
Set-TextValue $ws.Range('D35') '1.373'
Set-TextValue $ws.Range('E35') '  -4.00%  '
Set-TextValue $ws.Range('E36') '  -0.84%  '
Set-TextValue $ws.Range('D37') '0.02246'
Set-TextValue $ws.Range('E37') '  -0.65%  '
Set-TextValue $ws.Range('D38') '0.06040'
Set-TextValue $ws.Range('E38') '  -1.22%  '
Set-TextValue $ws.Range('D39') '1.170'
Set-TextValue $ws.Range('E39') '  -0.88%  '
Set-TextValue $ws.Range('D40') '7.921'
Set-TextValue $ws.Range('E40') '  -6.28%  '
Set-TextValue $ws.Range('D41') '0.5831'
Set-TextValue $ws.Range('E41') '  -1.04%  '
Set-TextValue $ws.Range('D42') '0.1833'
Set-TextValue $ws.Range('E42') '  -0.01%  '
Set-TextValue $ws.Range('D43') '10.03'
Set-TextValue $ws.Range('E43') '  -1.76%  '
Set-TextValue $ws.Range('D44') '1.285'
Set-TextValue $ws.Range('E44') '  +2.66%  '
Set-TextValue $ws.Range('B45') 'Cronos'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D45') '0.07701'
Set-TextValue $ws.Range('E45') '  +2.37%  '
Set-TextValue $ws.Range('B46') 'RenderToken'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D46') '2.361'
Set-TextValue $ws.Range('E46') '  -0.65%  '
Set-TextValue $ws.Range('D47') '12.21'
Set-TextValue $ws.Range('E47') '  +0.05%  '
Set-TextValue $ws.Range('D48') '0.5474'
Set-TextValue $ws.Range('E48') '  -1.81%  '
Set-TextValue $ws.Range('D49') '1.908'
Set-TextValue $ws.Range('E49') '  -1.31%  '
Set-TextValue $ws.Range('D50') '113.13'
Set-TextValue $ws.Range('E50') '  +0.17%  '
Set-TextValue $ws.Range('D51') '0.2939'
Set-TextValue $ws.Range('E51') '  +0.12%  '

$excel.CutCopyMode = $false
Write-Host "Done applying changes"